$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2020" data column has arrived for this indicator. Extend the table
# by one column (O), mirroring the existing "2019" column's (N) formatting
# for both the year header (row 4) and the data value (row 5).
$ws.Range("N4").Copy($ws.Range("O4")) | Out-Null
$ws.Range("O4").Value = 2020

$ws.Range("N5").Copy($ws.Range("O5")) | Out-Null
$ws.Range("O5").Value = 83.3

# Reflect the post-edit cursor position saved in the workbook.
$ws.Range("O12").Select() | Out-Null
